# abstract out detect_decel_onset and detect_brake_jerk_end for reusability
#
# - params sheet: split the generic JERK_NEG_THD / JERK_POS_THD parameters
#   into an AEB-specific negative-jerk threshold (+ a new latency-window
#   parameter) and rename the remaining pair to the BrakeJerk-specific
#   names, updating their descriptions accordingly.
# - makes "params" the active sheet/tab instead of "KPI".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("params")

# Insert two new rows right before the old "PRE_TIME_FCW" row (row 13),
# pushing everything below it down by two rows.
$ws.Rows("13:14").Insert() | Out-Null

# New row 13: AEB-specific negative jerk threshold.
$ws.Range("A13").Value = "AEB_JERK_NEG_THD"
$ws.Range("B13").Value = -30
$ws.Range("C13").Value = "float"
$ws.Range("D13").Value = "m/s³"
$ws.Range("E13").Value = "Negative jerk threshold to detect the start of AEB intervention "
$ws.Range("F13").Value = "AebKpiExtractor"

# New row 14: latency window (sample count) used for AEB latency detection.
$ws.Range("A14").Value = "LATENCY_WINDOW_SAMPLES"
$ws.Range("B14").Value = 30
$ws.Range("C14").Value = "int"
$ws.Range("E14").Value = "Sample window after AEB start for latency detection"
$ws.Range("F14").Value = "AebKpiExtractor"

# The old generic JERK_NEG_THD / JERK_POS_THD rows (now shifted to rows
# 18-19) become BrakeJerk-specific, with clarified descriptions.
$ws.Range("A18").Value = "BRAKEJERK_JERK_NEG_THD"
$ws.Range("E18").Value = "Negative jerk threshold to detect the start of BrakeJerk intervention "

$ws.Range("A19").Value = "BRAKEJERK_JERK_POS_THD"
$ws.Range("E19").Value = "Positive jerk threshold to detect the end of BrakeJerk intervention"

# Column A got wider to fit the new, longer parameter names.
$ws.Columns("A").ColumnWidth = 27.166666666666668

# Make "params" the active sheet/selection, replacing "KPI".
$ws.Activate() | Out-Null
$ws.Range("C26").Select() | Out-Null

Write-Output "done"
